# Apply the "最大图片大小" (max image size) config change:
#   2097152 bytes (2MB)  ->  77824 bytes (76KB)
# on sheet "公共参数|Common", row 22 (MAX_IMAGE_SIZE).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Client / server values for MAX_IMAGE_SIZE (columns C and D)
$ws.Range("C22").Value2 = 77824
$ws.Range("D22").Value2 = 77824

# Description text in column F (backed by the shared string table)
$ws.Range("F22").Value2 = "最大图片大小76KB"

# Move the active selection from F13 to F23, as recorded in the sheet view
$ws.Range("F23").Select()
